$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that can look like numbers (e.g. "146.39");
# force text format so Excel stores them verbatim instead of coercing to a float.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.028.73"
$ws.Range("E2").Value = "  -3.83%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.134.14"
$ws.Range("E3").Value = "  -3.57%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.38"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.39"
$ws.Range("E6").Value = "  -7.20%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.137.31"
$ws.Range("E8").Value = "  -3.50%  "
$ws.Range("E9").Value = "  -4.73%  "
$ws.Range("E10").Value = "  -6.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.48"
$ws.Range("E11").Value = "  -6.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.473"
$ws.Range("E12").Value = "  -6.17%  "
$ws.Range("E13").Value = "  -8.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.42"
$ws.Range("E14").Value = "  -9.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.663.24"
$ws.Range("E15").Value = "  -3.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.049.86"
$ws.Range("E16").Value = "  -3.98%  "
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.143.70"
$ws.Range("E18").Value = "  -3.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.89"
$ws.Range("E19").Value = "  -7.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "478.67"
$ws.Range("E20").Value = "  -5.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.79"
$ws.Range("E21").Value = "  -3.98%  "
$ws.Range("E22").Value = "  -5.33%  "
$ws.Range("E23").Value = "  -4.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.55"
$ws.Range("E24").Value = "  -8.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.30"
$ws.Range("E25").Value = "  -3.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("E27").Value = "  -5.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.39"
$ws.Range("E28").Value = "  -7.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.18"
$ws.Range("E29").Value = "  -9.75%  "
$ws.Range("E30").Value = "  -1.20%  "
$ws.Range("E31").Value = "  -19.01%  "
$ws.Range("B32").Value = "Stacks"
$ws.Range("C32").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.73"
$ws.Range("E32").Value = "  -6.27%  "
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("E34").Value = "  -7.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.10"
$ws.Range("E35").Value = "  -4.99%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.11"
$ws.Range("E36").Value = "  -2.81%  "
$ws.Range("E37").Value = "  -7.95%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "462.04"
$ws.Range("E38").Value = "  -6.46%  "
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0728"
$ws.Range("E39").Value = "  -8.48%  "
$ws.Range("E40").Value = "  -13.80%  "
$ws.Range("E41").Value = "  -7.43%  "
$ws.Range("E42").Value = "  -8.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.40"
$ws.Range("E43").Value = "  -4.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.840.00"
$ws.Range("E44").Value = "  -4.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.265"
$ws.Range("E45").Value = "  -9.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.25"
$ws.Range("E46").Value = "  -10.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.45"
$ws.Range("E47").Value = "  -8.27%  "
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("E49").Value = "  -7.18%  "
$ws.Range("E50").Value = "  -4.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "119.49"
$ws.Range("E51").Value = "  -1.60%  "
